# Fixed parts orientation in PnP file.
# Swap the Rotation (column E) values between the two groups of parts:
#  - rows 2-4 and 23 go from 180 -> 0
#  - rows 5-22 go from 0 -> 180

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2:E4").Value2 = 0
$ws.Range("E5:E22").Value2 = 180
$ws.Range("E23").Value2 = 0

# Update the selected cell to match the saved selection in the workbook.
$ws.Range("F5").Select() | Out-Null
